# "fixed levels that I screwed up"
# The level-design grid (rows 2-25, columns F/G/H = lane placements) had a
# batch of leftover "Flamingo" placeholders that needed to become "Gnome"
# (or be cleared out), plus several missing "Gnome" placements were added
# back into the block-1 timeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: the block-1 header row had "Flamingo" in all three lanes.
# F4 and G4 get cleared out entirely; H4 becomes "Gnome". G4 additionally
# picks up the same formatting already used by the other "marker" cells
# in this sheet (style index 12 in the original file) -- copy that
# formatting from a cell that already carries it (H2), then clear the
# copied content back out.
$ws.Range("H2").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").ClearContents() | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F4").ClearContents() | Out-Null
$ws.Range("H4").Value = "Gnome"

# --- Scattered "Gnome" placements added across the existing grid.
# These cells already exist (possibly blank) in the sheet, so writing the
# value in place preserves whatever formatting they already had.
$gnomeCells = @("H2", "G3", "G5", "F7", "F8", "G9", "G10", "H12", "G13", `
                "H14", "G15", "F17", "F18", "G19", "H22", "G23", "H24", "G25")
foreach ($addr in $gnomeCells) {
    $ws.Range($addr).Value = "Gnome"
}

# --- Row 20 previously stopped at column E; F20:H20 are brand-new cells.
# Give them plain/default formatting (matching the rest of the row) rather
# than whatever a freshly-materialized cell would otherwise inherit.
$ws.Range("F20").Value = ""
$ws.Range("F20").Style = "Normal"
$ws.Range("H20").Value = ""
$ws.Range("H20").Style = "Normal"
$ws.Range("G20").Value = "Gnome"
$ws.Range("G20").Style = "Normal"

# --- Selection moved to E20 in the saved file.
$ws.Range("E20").Select() | Out-Null
